$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 199 (old rows 199-257 shift down to 201-259)
$ws.Rows("199:200").Insert()

# New row 199: Candy White / Primera
$ws.Range("A199").Value = 11
$ws.Range("B199").Value = "Vega Monumental Concepción"
$ws.Range("C199").Value = "Bíobío"
$ws.Range("D199").Value = 44588
$ws.Range("D199").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E199").Value = 8
$ws.Range("F199").Value = "Fruta"
$ws.Range("G199").Value = 100103
$ws.Range("H199").Value = "Frutos de hueso (carozo)"
$ws.Range("I199").Value = 100103006
$ws.Range("J199").Value = "Nectarín"
$ws.Range("K199").Value = "Candy White"
$ws.Range("L199").Value = "Primera"
$ws.Range("M199").Value = 220
$ws.Range("N199").Value = 9000
$ws.Range("O199").Value = 9500
$ws.Range("P199").Value = 9227
$ws.Range("Q199").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R199").Value = "Región de O'Higgins"
$ws.Range("S199").Value = 577
$ws.Range("T199").Value = 16

# New row 200: Sun Rise / Primera
$ws.Range("A200").Value = 11
$ws.Range("B200").Value = "Vega Monumental Concepción"
$ws.Range("C200").Value = "Bíobío"
$ws.Range("D200").Value = 44588
$ws.Range("D200").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E200").Value = 8
$ws.Range("F200").Value = "Fruta"
$ws.Range("G200").Value = 100103
$ws.Range("H200").Value = "Frutos de hueso (carozo)"
$ws.Range("I200").Value = 100103006
$ws.Range("J200").Value = "Nectarín"
$ws.Range("K200").Value = "Sun Rise"
$ws.Range("L200").Value = "Primera"
$ws.Range("M200").Value = 250
$ws.Range("N200").Value = 9500
$ws.Range("O200").Value = 10000
$ws.Range("P200").Value = 9760
$ws.Range("Q200").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R200").Value = "Región de O'Higgins"
$ws.Range("S200").Value = 610
$ws.Range("T200").Value = 16
